$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 22-30 down/around: social influence block moves up one row, and
# "Empirical studies of simple contagion" moves to B23 (holiday week),
# pushing everything else to reflect the new order (holiday coming up).

# Row 22: becomes Social influence, herding, and cascades (was row 24's content)
$ws.Range("A22").Value = "[Social influence](#sec:socialinfluence)"
$ws.Range("D22").Value = "Social influence, herding, and cascades"
$ws.Range("G22").Value = "? Mini Project 04"

# Row 23: NO CLASS (HOLIDAY) stays, but B23 now holds the empirical studies text
$ws.Range("B23").Value = "Empirical studies of simple contagion"
$ws.Range("D23").Value = "NO CLASS (HOLIDAY)"

# Row 24: A24 cleared, D24 becomes Threshold models and complex contagion
$ws.Range("A24").ClearContents()
$ws.Range("D24").Value = "Threshold models and complex contagion"

# Row 25: A25 gets the Dynamics link, D25 becomes Complex contagion on networks
$ws.Range("A25").Value = "[Dynamics: Complex contagion and social influence](#sec:complexcontagion)"
$ws.Range("D25").Value = "Complex contagion on networks"

# Row 26: A26 cleared, D26 becomes Complex contagion on networks, cont.; G26 unchanged
$ws.Range("A26").ClearContents()
$ws.Range("D26").Value = "Complex contagion on networks, cont."

# Row 27 & 28 unchanged (NO CLASS / THANKSGIVING)

# Row 29: A29 gets the Cooperation link, D29 becomes Cooperation and networks
$ws.Range("A29").Value = "[Cooperation](#sec:cooperation)"
$ws.Range("D29").Value = "Cooperation and networks"

# Row 30: A30 cleared, D30 becomes Empirical studies / Wrap up
$ws.Range("A30").ClearContents()
$ws.Range("D30").Value = "Empirical studies / Wrap up"
